$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("membership")

# --- Header row: add two new date placeholder columns after "${currency}" ---
$ws.Range("F1").Value = '${date1}'
$ws.Range("G1").Value = '${date2}'

# Match the formatting of the existing header cells (style reuse, not a
# brand-new style) by copying E1's format onto the new header cells.
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1:G1").PasteSpecial(-4122) | Out-Null

# --- Sample data row: add the corresponding date values ---
$ws.Range("F2").Value2 = 44637.0
$ws.Range("G2").Value2 = 44637.0

# Start from the neighboring data cell's format (same font/alignment as the
# rest of the row), then apply the date-specific number format.
$ws.Range("D2").Copy() | Out-Null
$ws.Range("F2:G2").PasteSpecial(-4122) | Out-Null
$ws.Range("F2:G2").NumberFormat = "yyyy-mm-dd"

$excel.CutCopyMode = 0
